$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.928.06"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "2.670.47"
$ws.Range("E3").Value = "  +7.58%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'113.82"
$ws.Range("E5").Value = "  +8.64%  "
$ws.Range("D6").Value = "'326.46"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").Value = "'0.529"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").Value = "'41.08"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("D11").Value = "'20.08"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'0.0825"
$ws.Range("E12").Value = "  +3.23%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'7.37"
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("D15").Value = "3.095.06"
$ws.Range("E15").Value = "  +7.81%  "
$ws.Range("D16").Value = "2.678.37"
$ws.Range("E16").Value = "  +7.50%  "
$ws.Range("D17").Value = "'0.877"
$ws.Range("E17").Value = "  +6.58%  "
$ws.Range("D18").Value = "49.890.99"
$ws.Range("E18").Value = "  +4.25%  "
$ws.Range("D19").Value = "'13.19"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").Value = "'6.79"
$ws.Range("E20").Value = "  +4.10%  "
$ws.Range("D21").Value = "'2.89"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("D23").Value = "'72.51"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'277.48"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").Value = "'26.87"
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +5.83%  "
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "'36.24"
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").Value = "'50.28"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'19.68"
$ws.Range("E34").Value = "  +4.27%  "
$ws.Range("D35").Value = "'0.0818"
$ws.Range("E35").Value = "  +6.54%  "
$ws.Range("D36").Value = "'5.11"
$ws.Range("E36").Value = "  +13.18%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +7.40%  "
$ws.Range("D39").Value = "'3.16"
$ws.Range("E39").Value = "  +10.65%  "
$ws.Range("D40").Value = "'124.74"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.113"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'22.68"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +6.53%  "
$ws.Range("D45").Value = "2.117.20"
$ws.Range("E45").Value = "  +6.23%  "
$ws.Range("D46").Value = "'3.31"
$ws.Range("E46").Value = "  +5.37%  "
$ws.Range("D47").Value = "'2.26"
$ws.Range("E47").Value = "  +14.68%  "
$ws.Range("E48").Value = "  +7.41%  "
$ws.Range("D49").Value = "'9.05"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "'5.39"
$ws.Range("E50").Value = "  +5.13%  "
$ws.Range("D51").Value = "'59.61"
$ws.Range("E51").Value = "  +6.71%  "
